$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "Criticality indicators (6_CR) by material"
$ws.Range("H18").Value = "Criticality"
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = "material"
$ws.Range("K18").Value = 26

$ws.Range("M14").Select()
